# Clean up the raw CSV-import sheet: drop the numeric index / Address / Extra
# columns, and move Area Code to the front so the final layout is:
#   Area Code | First | Last | City | State | Income

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the values we still need before they get overwritten ---
$areaCode = @{}
$city     = @{}
$state    = @{}
$income   = @{}
for ($r = 2; $r -le 7; $r++) {
    $areaCode[$r] = $ws.Range("G$r").Value2   # column G
    $city[$r]     = $ws.Range("E$r").Value2   # column E
    $state[$r]    = $ws.Range("F$r").Value2   # column F
    $income[$r]   = $ws.Range("H$r").Value2   # column H
}

# --- header row ---
$ws.Range("A1").Value = "Area Code"
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats - pick up the bold/centered header style

$ws.Range("D1").Value = "City"
$ws.Range("E1").Value = "State"
$ws.Range("F1").Value = "Income"

# --- data rows ---
for ($r = 2; $r -le 7; $r++) {
    $ws.Range("A$r").Value = $areaCode[$r]
    $ws.Range("D$r").Value = $city[$r]
    $ws.Range("E$r").Value = $state[$r]
    $ws.Range("F$r").Value = $income[$r]
}

# fix up the names that previously carried embedded quotes/commas
$ws.Range("B4").Value = "John"
$ws.Range("B6").Value = "N/A"
$ws.Range("B7").Value = "Joan"

# --- drop the now-unused trailing columns (old Area Code / Extra) ---
$ws.Range("G1:H7").Clear()

$ws.Range("A1").Select()
